# Add more county rows to the "County" sheet so there's enough data to
# compute bins, per the commit message ("Add enough data to the insecurity
# test data that we can bin it").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("County")

# New counties to append below the existing "Jo Daviess County, Illinois" row.
# Columns: A = County, State ; B = 2018 Food Insecurity % ;
#          C = 2020 (revised) Food Insecurity % ; D = 2018 Child Food Insecurity % ;
#          E = 2020 (revised) Child Food Insecurity %
$rows = @(
    @("Johnson County, Illinois",   0.119, 0.163, 0.182, 0.268),
    @("Kane County, Illinois",      0.071, 0.116, 0.094, 0.183),
    @("Kankakee County, Illinois",  0.116, 0.157, 0.157, 0.238),
    @("Kendall County, Illinois",   0.048, 0.092, 0.064, 0.15)
)

$startRow = 3
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Range("A$r").Value = $data[0]
    $ws.Range("B$r").Value = $data[1]
    $ws.Range("C$r").Value = $data[2]
    $ws.Range("D$r").Value = $data[3]
    $ws.Range("E$r").Value = $data[4]

    # Match the number format already used in the data column (D uses the
    # three-decimal custom format; the rest stay General like row 2).
    $ws.Range("D$r").NumberFormat = "0.000"
}

# Keep the same active-cell convention as the source edit: selection moves
# just past the newly added data.
$ws.Range("A11").Select()
